$d = $word.ActiveDocument
$end = $d.Content.End
$r = $d.Range($end, $end)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">---</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Store Builder &amp; Customization Expansion</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Updated: 2026-02-18</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"></w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Module Name</w:t><w:tab/><w:t xml:space="preserve">Developed</w:t><w:tab/><w:t xml:space="preserve">Partial Developed</w:t><w:tab/><w:t xml:space="preserve">Need To Develop</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Theme / Design Engine</w:t><w:tab/><w:t xml:space="preserve">Theme marketplace API + platform owner theme create/list UI + store theme apply with plan gating; live preview runtime route `/s/:subdomain/*`; section-based homepage layout persistence + editor; custom JSON token controls (header/footer/banner/design)</w:t><w:tab/><w:t xml:space="preserve">Visual drag-drop is reorder-based (up/down + section editor), advanced WYSIWYG widget tooling is partial</w:t><w:tab/><w:t xml:space="preserve">Full visual drag canvas with nested sections, versioned theme publishing workflow, safe sandbox for merchant custom JS</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Branding Controls</w:t><w:tab/><w:t xml:space="preserve">Logo/favicon upload API (`/storefront/media/upload`) + media asset registry, colors/typography via design tokens JSON, banners/sliders via JSON blocks, storefront runtime uses active branding</w:t><w:tab/><w:t xml:space="preserve">Media pipeline currently local static storage (or ASSET_BASE_URL mapping) and not fully external CDN-managed</w:t><w:tab/><w:t xml:space="preserve">Managed CDN integration, asset optimization/transform pipeline, advanced typography presets/UI controls</w:t></w:r></w:p><w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Storefront Runtime + B2B Visibility</w:t><w:tab/><w:t xml:space="preserve">Public storefront APIs + frontend renderer consume theme/layout/navigation/pages/products; pricing visibility flags wired (`showPricing`, `loginToViewPrice`, `catalogMode`, `catalogVisibilityJson`)</w:t><w:tab/><w:t xml:space="preserve">Customer-specific catalog visibility still JSON-driven (no business-rule UI), login-gated pricing policy needs end-user auth enforcement on public storefront</w:t><w:tab/><w:t xml:space="preserve">Complete B2B rule engine + customer-group catalog filters + enforced authenticated B2B storefront sessions</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
